$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from H1 into the new I1/J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF)
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 5
